$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows 23-26 mirror the formatting of row 22 (same column styles: A/B/D/E plain
# bordered cells, C bordered cell). Copy the format down first, then fill in values.
[void]$ws.Range("A22:E22").Copy()
[void]$ws.Range("A23:E26").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A23").Value = "Notifications022"
$ws.Range("B23").Value = "OPQA-1102"
$ws.Range("C23").Value = "Verify that system is able to recommend six people for user"
$ws.Range("D23").Value = "N"
$ws.Range("E23").Value = "SKIP"

$ws.Range("A24").Value = "Notifications023"
$ws.Range("B24").Value = "OPQA-211"
$ws.Range("C24").Value = "Verify that user is able to view top commenters information in home page"
$ws.Range("D24").Value = "N"
$ws.Range("E24").Value = "SKIP"

$ws.Range("A25").Value = "Notifications024"
$ws.Range("B25").Value = "OPQA-212"
$ws.Range("C25").Value = "Verify that user is able to view Most viewed documents in home page"
$ws.Range("D25").Value = "N"
$ws.Range("E25").Value = "SKIP"

$ws.Range("A26").Value = "Notifications025"
$ws.Range("B26").Value = "OPQA-1103"
$ws.Range("C26").Value = "Verify that system is able to recommend three articles for user"
$ws.Range("D26").Value = "N"
$ws.Range("E26").Value = "SKIP"

# Update the saved view state to match: active cell moved down one row to line up
# with the newly appended rows.
[void]$ws.Range("D6").Select()
